$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so number-like strings
# ("281.87", "0.9671", etc.) are not reinterpreted as numeric values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '21.133.46'
$ws.Range("E2").Value = '  +3.54%  '
$ws.Range("D3").Value = '1.536.76'
$ws.Range("E3").Value = '  +5.09%  '
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '0.9661'
$ws.Range("E5").Value = '  +1.52%  '
$ws.Range("D6").Value = '281.87'
$ws.Range("E6").Value = '  +2.51%  '
$ws.Range("D7").Value = '0.3622'
$ws.Range("E7").Value = '  -0.95%  '
$ws.Range("D8").Value = '0.3168'
$ws.Range("E8").Value = '  +3.39%  '
$ws.Range("D9").Value = '40.61'
$ws.Range("E9").Value = '  +1.99%  '
$ws.Range("E10").Value = '  +5.60%  '
$ws.Range("E11").Value = '  +3.41%  '
$ws.Range("D12").Value = '1.005'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").Value = '5.657'
$ws.Range("E13").Value = '  +4.35%  '
$ws.Range("D14").Value = '18.69'
$ws.Range("E14").Value = '  +3.52%  '
$ws.Range("D15").Value = '6.341'
$ws.Range("E15").Value = '  +3.24%  '
$ws.Range("B16").Value = 'Dai'
$ws.Range("C16").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D16").Value = '0.9666'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001043'
$ws.Range("E17").Value = '  +1.75%  '
$ws.Range("D18").Value = '1.529.58'
$ws.Range("E18").Value = '  +4.54%  '
$ws.Range("D19").Value = '0.06067'
$ws.Range("E19").Value = '  +3.84%  '
$ws.Range("D20").Value = '72.28'
$ws.Range("E20").Value = '  +4.50%  '
$ws.Range("D21").Value = '5.696'
$ws.Range("E21").Value = '  +4.88%  '
$ws.Range("E22").Value = '  +3.96%  '
$ws.Range("D23").Value = '11.35'
$ws.Range("E23").Value = '  +4.30%  '
$ws.Range("D24").Value = '2.327'
$ws.Range("E24").Value = '  +3.82%  '
$ws.Range("D25").Value = '21.167.57'
$ws.Range("E25").Value = '  +3.64%  '
$ws.Range("D26").Value = '148.28'
$ws.Range("E26").Value = '  +4.37%  '
$ws.Range("D27").Value = '2.212'
$ws.Range("E27").Value = '  +6.63%  '
$ws.Range("D28").Value = '17.62'
$ws.Range("E28").Value = '  +2.87%  '
$ws.Range("D29").Value = '1.695.52'
$ws.Range("E29").Value = '  +4.87%  '
$ws.Range("D30").Value = '118.47'
$ws.Range("E30").Value = '  +4.77%  '
$ws.Range("D31").Value = '4.007'
$ws.Range("E31").Value = '  +4.08%  '
$ws.Range("D32").Value = '0.8488'
$ws.Range("E32").Value = '  +7.14%  '
$ws.Range("D33").Value = '5.171'
$ws.Range("E33").Value = '  +5.54%  '
$ws.Range("D34").Value = '0.07990'
$ws.Range("E34").Value = '  +1.34%  '
$ws.Range("D35").Value = '1.506'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = '1.208'
$ws.Range("E36").Value = '  +5.29%  '
$ws.Range("D37").Value = '4.940'
$ws.Range("E37").Value = '  +5.54%  '
$ws.Range("D38").Value = '0.05874'
$ws.Range("E38").Value = '  +2.88%  '
$ws.Range("D39").Value = '0.02099'
$ws.Range("E39").Value = '  +3.43%  '
$ws.Range("D40").Value = '10.62'
$ws.Range("E40").Value = '  +2.89%  '
$ws.Range("D41").Value = '7.718'
$ws.Range("E41").Value = '  +2.84%  '
$ws.Range("D42").Value = '0.9671'
$ws.Range("E42").Value = '  +0.92%  '
$ws.Range("D43").Value = '0.1911'
$ws.Range("E43").Value = '  +2.84%  '
$ws.Range("D44").Value = '0.5431'
$ws.Range("E44").Value = '  +3.26%  '
$ws.Range("D45").Value = '12.52'
$ws.Range("E45").Value = '  +5.29%  '
$ws.Range("D46").Value = '3.568'
$ws.Range("E46").Value = '  +2.35%  '
$ws.Range("D47").Value = '0.5429'
$ws.Range("E47").Value = '  +5.35%  '
$ws.Range("D48").Value = '121.28'
$ws.Range("E48").Value = '  +3.31%  '
$ws.Range("E49").Value = '  +6.88%  '
$ws.Range("D50").Value = '0.06573'
$ws.Range("E50").Value = '  +2.37%  '
$ws.Range("D51").Value = '0.9918'
$ws.Range("E51").Value = '  +0.09%  '

# Restore the default (unstyled) cell style on column D now that the
# values are committed as text, matching the original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"
